$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump the "Date" value (row 8, col B)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2024-08-22T16:27:20+00:00"

# ---------------------------------------------------------------------------
# 2. Elements sheet: row 6 ("Extension.value[x]") picks up a type-slicing
#    discriminator + rules, and its Type(s)/Short text change; then a brand
#    new row 7 is added for the "valueQuantity" slice of Extension.value[x].
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# -- 2a. Update row 6 in place --------------------------------------------
$ws.Cells.Item(6, 11).Value = "date`nQuantity"   # K6  Type(s)
$ws.Cells.Item(6, 12).Value = "Value of extension" # L6  Short
$ws.Cells.Item(6, 28).Value = "type:`$this}`n"    # AB6 Slicing Discriminator
$ws.Cells.Item(6, 29).Value = ""                  # AC6 Slicing Description (cleared)
$ws.Cells.Item(6, 31).Value = "open"              # AE6 Slicing Rules

# Re-apply row 6's original cell formatting (AB6/AC6 had no content before,
# so writing into them resets their style/number-format; pulling the format
# back from a still-untouched cell on the same row keeps everything uniform).
$ws.Range("A6:AK6").Copy()
$ws.Range("A6:AK6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -- 2b. Populate the new row 7 --------------------------------------------
$row7 = @{
  1  = "Extension.value[x]:valueQuantity"
  2  = "Extension.value[x]"
  3  = "valueQuantity"
  6  = "0"
  7  = "1"
  11 = "Quantity`n"
  12 = "Indicate age via relative date time extension or official date of when last vital status was assessed."
  13 = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4B/extensibility.html) for a list)."
  32 = "Extension.value[x]"
  33 = "0"
  34 = "1"
  36 = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`n"
  37 = "N/A"
}

# Columns F/G/AG/AH hold text that *looks* numeric ("0"/"1"); force the cell
# to text first so the engine keeps them as strings instead of numbers.
foreach ($col in 6, 7, 33, 34) {
  $ws.Cells.Item(7, $col).NumberFormat = "@"
}

foreach ($col in $row7.Keys) {
  $ws.Cells.Item(7, $col).Value = $row7[$col]
}

# Give every cell in the new row the same look (border/wrap/alignment/style)
# as the rest of the data rows by cloning row 6's formatting down onto it.
$ws.Range("A6:AK6").Copy()
$ws.Range("A7:AK7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The wrapped multi-line text we just dropped into row 6/7 makes the engine
# stamp an explicit row height; snap both rows back to the sheet's natural
# (default) height so no stray customHeight sticks around.
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()

# -- 2c. Column width tweaks (ID column widens, Slice Name column widens) --
$ws.Columns.Item(1).ColumnWidth = 31.7
$ws.Columns.Item(3).ColumnWidth = 14
